{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The title paragraph is the first paragraph in the document:\n// \"Python-Numpy Assignment 1\" -> \"Numpy Assignment 2\"\nconst titlePara = paragraphs.items[0];\nconst titleRange = titlePara.getRange();\ntitleRange.load(\"text\");\nawait context.sync();\n\nif (titleRange.text.indexOf(\"Python-Numpy Assignment 1\") !== -1) {\n  titleRange.insertText(\"Numpy Assignment 2\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Title paragraph: \"Python-Numpy Assignment 1\" -> \"Numpy Assignment 2\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Python-Numpy Assignment 1\"\n$find.Replacement.Text = \"Numpy Assignment 2\"\n$find.Execute(\n    [ref]\"Python-Numpy Assignment 1\",\n    [ref]$false,\n    [ref]$true,\n    [ref]$false,\n    [ref]$false,\n    [ref]$false,\n    [ref]$true,\n    [ref]1,\n    [ref]$false,\n    [ref]\"Numpy Assignment 2\",\n    [ref]2\n) | Out-Null\n"}
